$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# KiCost re-ran the export: the "Prj date:" (human-readable project timestamp)
# and "$ date:" (currency-rate lookup timestamp) strings were refreshed.
$ws.Range("B3").Value = "Sat Jun  4 23:26:00 2022"
$ws.Range("B4").Value = "2022-06-04 23:26:02"
